$wb = $excel.ActiveWorkbook

# --- "Instructions" sheet: scroll the view so column F (instead of F31) is the
#     top-left visible column again (view was left scrolled down to F31).
$wsInstructions = $wb.Worksheets.Item("Instructions")
$winInstructions = $wsInstructions.Application.Windows.Item(1)
$winInstructions.ScrollColumn = 6
$winInstructions.ScrollRow = 1

# --- "T-states (NEEDS UPDATE)" sheet: slightly taller header row.
$wsTStates = $wb.Worksheets.Item("T-states (NEEDS UPDATE)")
$wsTStates.Rows.Item(1).RowHeight = 47.4

# --- "Notes" sheet: add a new "Stack" note (with its address range) just above
#     the existing "If ROM boot is enabled:" block, pushing everything below it
#     down by two rows.
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Rows("10:11").Insert()
$wsNotes.Range("A9").Value = "Stack"
$wsNotes.Range("B10").Value = '$FF00-$FFFF'
$wsNotes.Range("B10").Select()
